$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.022.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.242.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.21%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "99.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +19.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "270.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.94%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.53%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.642"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0945"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +17.78%  "
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.562.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.824"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.248.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.997.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("E19").Value = "  +2.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.24%  "
$ws.Range("E22").Value = "  -2.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.56%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +12.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.65%  "
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0933"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.54%  "
$ws.Range("E33").Value = "  +3.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.23%  "
$ws.Range("E35").Value = "  +1.57%  "
$ws.Range("E36").Value = "  +1.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0352"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +34.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.247"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +24.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.71%  "
$ws.Range("E42").Value = "  +4.32%  "
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("E45").Value = "  +4.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.94%  "
$ws.Range("E48").Value = "  +4.60%  "
$ws.Range("E49").Value = "  +0.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.440"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("E51").Value = "  +1.15%  "
